$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-04-05"

# Update the label for the April row to reflect the new "through" date
$ws.Range("A5").Value = "April (through 04-05)"

# Update March 2022 value
$ws.Range("I4").Value = 134

# Update April row values for 2015, 2017, 2021, 2022
$ws.Range("B5").Value = 3
$ws.Range("D5").Value = 6
$ws.Range("H5").Value = 12
$ws.Range("I5").Value = 16

# Update Total row values for 2015, 2017, 2021, 2022
$ws.Range("B6").Value = 69
$ws.Range("D6").Value = 195
$ws.Range("H6").Value = 435
$ws.Range("I6").Value = 450
